# Weekly update: insert a new price-report row for "Achicoria" at
# Vega Central Mapocho de Santiago, which shifts every subsequent
# row (old 34..77) down by one (new 35..78) and grows the used range
# from A1:R77 to A1:R78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 34..77 down one slot, leaving a blank row 34 to fill in.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with this week's record (same
# market/category/quality pattern as its neighbours, new date + volume).
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 45128
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 100112010
$ws.Range("G34").Value = "Achicoria"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 90
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = 7000
$ws.Range("N34").Value = "$/caja 16 unidades"
$ws.Range("O34").Value = "Provincia de Quillota"
$ws.Range("P34").Value = 438
$ws.Range("Q34").Value = 16
$ws.Range("R34").Value = "Hortaliza"
